# Edit script: adds a "Top Losers" worksheet right after "Top Gainers",
# populates it with stock data, styles the header row to match the
# look of the other sheets, and bumps the "Last Updated" timestamp on
# the Metadata sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Top Losers" worksheet right after "Top Gainers"
#    so the final sheet order is:
#      Metadata, Top Gainers, Top Losers, 1 Month Performance,
#      distance from Dma50
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("Top Gainers")
$newSheet = $wb.Worksheets.Add($null, $afterSheet, 1)
$newSheet.Name = "Top Losers"

# Match the default page margins used by the other sheets in this workbook
# (0.75in / 0.75in / 1in / 1in / 0.5in / 0.5in), expressed in points.
$ps = $newSheet.PageSetup
$ps.LeftMargin = 0.75 * 72
$ps.RightMargin = 0.75 * 72
$ps.TopMargin = 1 * 72
$ps.BottomMargin = 1 * 72
$ps.HeaderMargin = 0.5 * 72
$ps.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------------
# 2. Clone the header formatting (bold font, thin borders, centered /
#    top-aligned) from the "Top Gainers" header row so the new sheet
#    re-uses the very same cell style instead of creating a new one.
# ---------------------------------------------------------------------
$srcHeader = $afterSheet.Range("A1:E1")
$srcHeader.Copy()
$dstHeader = $newSheet.Range("A1:E1")
$dstHeader.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Write the header row text.
# ---------------------------------------------------------------------
$headers = New-Object 'object[,]' 1,5
$headers[0,0] = "Icon"
$headers[0,1] = "Stock"
$headers[0,2] = "Latest"
$headers[0,3] = "Weekly"
$headers[0,4] = "Monthly"
$newSheet.Range("A1:E1").Value = $headers

# ---------------------------------------------------------------------
# 4. Write the 75 data rows (rows 2-76) in one bulk assignment.
# ---------------------------------------------------------------------
$data = New-Object 'object[,]' 75,5
$data[0,0] = "📉"
$data[0,1] = "IIFLCAPS"
$data[0,2] = -9.721500000000001
$data[0,3] = -5.9431
$data[0,4] = 20.4574
$data[1,0] = "📉"
$data[1,1] = "COHANCE"
$data[1,2] = -8.5204
$data[1,3] = -10.2636
$data[1,4] = -10.5392
$data[2,0] = "📉"
$data[2,1] = "KHAICHEM"
$data[2,2] = -8.4635
$data[2,3] = -6.7205
$data[2,4] = -4.6455
$data[3,0] = "📉"
$data[3,1] = "MOTILALOFS"
$data[3,2] = -7.6514
$data[3,3] = -2.7408
$data[3,4] = 12.7798
$data[4,0] = "📉"
$data[4,1] = "VERANDA"
$data[4,2] = -7.1035
$data[4,3] = -7.4096
$data[4,4] = 8.066700000000001
$data[5,0] = "📉"
$data[5,1] = "MOLDTKPAC"
$data[5,2] = -7.0505
$data[5,3] = -1.2462
$data[5,4] = 0.466
$data[6,0] = "📉"
$data[6,1] = "NETWEB"
$data[6,2] = -6.7177
$data[6,3] = 5.5715
$data[6,4] = 7.5565
$data[7,0] = "📉"
$data[7,1] = "FABTECH"
$data[7,2] = -6.147
$data[7,3] = 15.6232
$data[7,4] = "N/A"
$data[8,0] = "📉"
$data[8,1] = "CARTRADE"
$data[8,2] = -5.7253
$data[8,3] = 16.1193
$data[8,4] = 20.6359
$data[9,0] = "📉"
$data[9,1] = "TVSELECT"
$data[9,2] = -5.6153
$data[9,3] = -0.9738
$data[9,4] = -2.9968
$data[10,0] = "📉"
$data[10,1] = "NSLNISP"
$data[10,2] = -5.4542
$data[10,3] = 1.3037
$data[10,4] = 0.4681
$data[11,0] = "📉"
$data[11,1] = "NAM-INDIA"
$data[11,2] = -5.1992
$data[11,3] = -7.0279
$data[11,4] = -1.2034
$data[12,0] = "📉"
$data[12,1] = "KICL"
$data[12,2] = -5.0437
$data[12,3] = -0.8405
$data[12,4] = 21.7969
$data[13,0] = "📉"
$data[13,1] = "CCCL"
$data[13,2] = -5.0146
$data[13,3] = -4.576
$data[13,4] = -12.4759
$data[14,0] = "📉"
$data[14,1] = "CREDITACC"
$data[14,2] = -4.9692
$data[14,3] = -1.3216
$data[14,4] = 3.7319
$data[15,0] = "📉"
$data[15,1] = "KALAMANDIR"
$data[15,2] = -4.8415
$data[15,3] = 1.7451
$data[15,4] = 25.9996
$data[16,0] = "📉"
$data[16,1] = "CRAMC"
$data[16,2] = -4.7668
$data[16,3] = 5.978
$data[16,4] = "N/A"
$data[17,0] = "📉"
$data[17,1] = "SMLISUZU"
$data[17,2] = -4.7654
$data[17,3] = 4.993
$data[17,4] = -2.8236
$data[18,0] = "📉"
$data[18,1] = "MANAKCOAT"
$data[18,2] = -4.5883
$data[18,3] = -6.011
$data[18,4] = 24.8093
$data[19,0] = "📉"
$data[19,1] = "HDFCAMC"
$data[19,2] = -4.401
$data[19,3] = -2.6247
$data[19,4] = -2.4311
$data[20,0] = "📉"
$data[20,1] = "ATHERENERG"
$data[20,2] = -4.0945
$data[20,3] = -0.0142
$data[20,4] = 24.8806
$data[21,0] = "📉"
$data[21,1] = "SHAREINDIA"
$data[21,2] = -4.0806
$data[21,3] = -1.6889
$data[21,4] = 54.7217
$data[22,0] = "📉"
$data[22,1] = "INDIQUBE"
$data[22,2] = -4.0805
$data[22,3] = -4.8408
$data[22,4] = -6.7982
$data[23,0] = "📉"
$data[23,1] = "CHENNPETRO"
$data[23,2] = -3.973
$data[23,3] = 4.8359
$data[23,4] = 6.8488
$data[24,0] = "📉"
$data[24,1] = "KFINTECH"
$data[24,2] = -3.873
$data[24,3] = -1.7487
$data[24,4] = 7.3785
$data[25,0] = "📉"
$data[25,1] = "360ONE"
$data[25,2] = -3.8488
$data[25,3] = -4.976
$data[25,4] = 10.0293
$data[26,0] = "📉"
$data[26,1] = "SMSPHARMA"
$data[26,2] = -3.7339
$data[26,3] = -3.0871
$data[26,4] = 17.4387
$data[27,0] = "📉"
$data[27,1] = "BHARATWIRE"
$data[27,2] = -3.5327
$data[27,3] = 22.8336
$data[27,4] = 23.8979
$data[28,0] = "📉"
$data[28,1] = "ABSLAMC"
$data[28,2] = -3.5313
$data[28,3] = -5.9355
$data[28,4] = -1.2887
$data[29,0] = "📉"
$data[29,1] = "SUMMITSEC"
$data[29,2] = -3.4113
$data[29,3] = -1.6476
$data[29,4] = 6.0097
$data[30,0] = "📉"
$data[30,1] = "SPLPETRO"
$data[30,2] = -3.3984
$data[30,3] = -5.0241
$data[30,4] = -7.769
$data[31,0] = "📉"
$data[31,1] = "MPSLTD"
$data[31,2] = -3.2782
$data[31,3] = -4.6314
$data[31,4] = 2.1755
$data[32,0] = "📉"
$data[32,1] = "CAMS"
$data[32,2] = -3.2545
$data[32,3] = -0.6366000000000001
$data[32,4] = 2.5781
$data[33,0] = "📉"
$data[33,1] = "SPARC"
$data[33,2] = -3.1709
$data[33,3] = 4.8337
$data[33,4] = 6.3311
$data[34,0] = "📉"
$data[34,1] = "PRUDENT"
$data[34,2] = -3.127
$data[34,3] = -3.5103
$data[34,4] = 2.1213
$data[35,0] = "📉"
$data[35,1] = "ANANDRATHI"
$data[35,2] = -3.0775
$data[35,3] = -0.8672
$data[35,4] = 9.1835
$data[36,0] = "📉"
$data[36,1] = "NLCINDIA"
$data[36,2] = -3.0757
$data[36,3] = -4.5618
$data[36,4] = -11.6431
$data[37,0] = "📉"
$data[37,1] = "YATRA"
$data[37,2] = -3.0403
$data[37,3] = -2.8455
$data[37,4] = 7.3711
$data[38,0] = "📉"
$data[38,1] = "BOSCHLTD"
$data[38,2] = -3.0099
$data[38,3] = -3.123
$data[38,4] = -2.0055
$data[39,0] = "📉"
$data[39,1] = "DRREDDY"
$data[39,2] = -2.9859
$data[39,3] = -2.5475
$data[39,4] = 2.2228
$data[40,0] = "📉"
$data[40,1] = "ROSSTECH"
$data[40,2] = -2.9778
$data[40,3] = 1.9028
$data[40,4] = -6.8057
$data[41,0] = "📉"
$data[41,1] = "OAL"
$data[41,2] = -2.9496
$data[41,3] = -1.278
$data[41,4] = 8.7362
$data[42,0] = "📉"
$data[42,1] = "ENDURANCE"
$data[42,2] = -2.939
$data[42,3] = -2.2945
$data[42,4] = 3.4531
$data[43,0] = "📉"
$data[43,1] = "POLICYBZR"
$data[43,2] = -2.907
$data[43,3] = 2.2365
$data[43,4] = 1.2573
$data[44,0] = "📉"
$data[44,1] = "AYMSYNTEX"
$data[44,2] = -2.9052
$data[44,3] = -0.3705
$data[44,4] = -10.494
$data[45,0] = "📉"
$data[45,1] = "DIGITIDE"
$data[45,2] = -2.8795
$data[45,3] = 3.2317
$data[45,4] = 6.2968
$data[46,0] = "📉"
$data[46,1] = "RUBICON"
$data[46,2] = -2.8712
$data[46,3] = -2.9654
$data[46,4] = "N/A"
$data[47,0] = "📉"
$data[47,1] = "STARHEALTH"
$data[47,2] = -2.8707
$data[47,3] = -1.5572
$data[47,4] = 7.5434
$data[48,0] = "📉"
$data[48,1] = "KIRIINDUS"
$data[48,2] = -2.8411
$data[48,3] = -1.3849
$data[48,4] = 1.4335
$data[49,0] = "📉"
$data[49,1] = "UNIMECH"
$data[49,2] = -2.8008
$data[49,3] = -1.6104
$data[49,4] = -0.4585
$data[50,0] = "📉"
$data[50,1] = "TTKPRESTIG"
$data[50,2] = -2.7438
$data[50,3] = 8.001200000000001
$data[50,4] = 9.650499999999999
$data[51,0] = "📉"
$data[51,1] = "PFOCUS"
$data[51,2] = -2.7039
$data[51,3] = -2.6276
$data[51,4] = -1.2163
$data[52,0] = "📉"
$data[52,1] = "ALLDIGI"
$data[52,2] = -2.6342
$data[52,3] = -0.2306
$data[52,4] = -5.3103
$data[53,0] = "📉"
$data[53,1] = "PRIVISCL"
$data[53,2] = -2.6288
$data[53,3] = -2.1048
$data[53,4] = 19.7451
$data[54,0] = "📉"
$data[54,1] = "CANHLIFE"
$data[54,2] = -2.6148
$data[54,3] = 5.2953
$data[54,4] = "N/A"
$data[55,0] = "📉"
$data[55,1] = "GKENERGY"
$data[55,2] = -2.6122
$data[55,3] = -9.807700000000001
$data[55,4] = 23.2758
$data[56,0] = "📉"
$data[56,1] = "SGFIN"
$data[56,2] = -2.592
$data[56,3] = -0.06270000000000001
$data[56,4] = 11.7235
$data[57,0] = "📉"
$data[57,1] = "ARVINDFASN"
$data[57,2] = -2.549
$data[57,3] = -2.9892
$data[57,4] = -4.4223
$data[58,0] = "📉"
$data[58,1] = "EDELWEISS"
$data[58,2] = -2.5422
$data[58,3] = -3.3745
$data[58,4] = 8.5305
$data[59,0] = "📉"
$data[59,1] = "SAMHI"
$data[59,2] = -2.5284
$data[59,3] = 1.8231
$data[59,4] = 2.8516
$data[60,0] = "📉"
$data[60,1] = "TBOTEK"
$data[60,2] = -2.524
$data[60,3] = -3.5732
$data[60,4] = 1.036
$data[61,0] = "📉"
$data[61,1] = "UJJIVANSFB"
$data[61,2] = -2.5201
$data[61,3] = 0.3845
$data[61,4] = 12.6645
$data[62,0] = "📉"
$data[62,1] = "AMBER"
$data[62,2] = -2.5098
$data[62,3] = -0.1082
$data[62,4] = 2.763
$data[63,0] = "📉"
$data[63,1] = "GRPLTD"
$data[63,2] = -2.4898
$data[63,3] = -5.9894
$data[63,4] = -5.4586
$data[64,0] = "📉"
$data[64,1] = "NESCO"
$data[64,2] = -2.4722
$data[64,3] = 1.9934
$data[64,4] = 3.8931
$data[65,0] = "📉"
$data[65,1] = "PILANIINVS"
$data[65,2] = -2.4546
$data[65,3] = -0.7907
$data[65,4] = 4.267
$data[66,0] = "📉"
$data[66,1] = "NSIL"
$data[66,2] = -2.4088
$data[66,3] = -1.7646
$data[66,4] = 4.7431
$data[67,0] = "📉"
$data[67,1] = "COALINDIA"
$data[67,2] = -2.4016
$data[67,3] = -3.058
$data[67,4] = -2.0387
$data[68,0] = "📉"
$data[68,1] = "JNKINDIA"
$data[68,2] = -2.3482
$data[68,3] = -2.8371
$data[68,4] = 4.2622
$data[69,0] = "📉"
$data[69,1] = "FCL"
$data[69,2] = -2.3453
$data[69,3] = -2.616
$data[69,4] = -0.02
$data[70,0] = "📉"
$data[70,1] = "DEEDEV"
$data[70,2] = -2.3334
$data[70,3] = -6.6528
$data[70,4] = -7.4227
$data[71,0] = "📉"
$data[71,1] = "WEALTH"
$data[71,2] = -2.2793
$data[71,3] = -3.8356
$data[71,4] = -2.7981
$data[72,0] = "📉"
$data[72,1] = "RATNAMANI"
$data[72,2] = -2.2788
$data[72,3] = -0.4626
$data[72,4] = 0.8712
$data[73,0] = "📉"
$data[73,1] = "CSBBANK"
$data[73,2] = -2.2695
$data[73,3] = 2.3137
$data[73,4] = 10.6999
$data[74,0] = "📉"
$data[74,1] = "BBOX"
$data[74,2] = -2.2639
$data[74,3] = -4.7636
$data[74,4] = 5.1528
$newSheet.Range("A2:E76").Value = $data

# ---------------------------------------------------------------------
# 5. Update the "Last Updated" timestamp on the Metadata sheet.
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("A2").Value = "29 Oct 2025, 07:24 PM"
